# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml  -> bound (via the presentation + slide master
#                             relationships) to every slide that's on screen.
#                             It currently carries the "Integral" / "Red
#                             Violet" palette.
#   ppt/theme/theme2.xml  -> bound only to the notes master, currently
#                             carrying the stock "Office Theme" / "Office"
#                             palette.
#
# The authored edit swaps the two palettes: the slide-facing theme becomes
# the stock Office palette, and the notes-only theme becomes the old
# Integral palette. Apply the reachable (slide-facing) half of that swap
# through the Office theme-colors API: ThemeColorScheme exposes all twelve
# scheme slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) in a fixed
# order, and editing them through any slide rewrites the shared theme part
# used by the whole deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# PowerPoint's RGB color longs are packed 0x00BBGGRR (little-endian), i.e.
# value = R + G*256 + B*65536, matching the &H00bbggrr the Color dialog uses.
$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      44546A
$tcs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  4472C4
$tcs.Item(10).RGB = 4697456  # accent6  70AD47
$tcs.Item(11).RGB = 12673797 # hlink    0563C1
$tcs.Item(12).RGB = 7491477  # folHlink 954F72
